# Weekly fruit/vegetable price update:
#  - insert a new week's record at row 53 (pushing all subsequent rows down)
#  - insert a second new week's record near the end (at row 80, after the
#    first shift), pushing the last two existing rows down as well.
#
# Final sheet grows from 80 data+header rows (A1:R80) to 82 (A1:R82).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 53 -------------------------------------------------
$ws.Rows("53:53").Insert()

$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 44567
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = 100112001
$ws.Range("G53").Value = "Berenjena"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 300
$ws.Range("K53").Value = 9000
$ws.Range("L53").Value = 9000
$ws.Range("M53").Value = 9000
$ws.Range("N53").Value = "$/caja 50 unidades"
$ws.Range("O53").Value = "Región del Maule"
$ws.Range("P53").Value = 180
$ws.Range("Q53").Value = 50
$ws.Range("R53").Value = "Hortaliza"

# --- Insert new row 80 (post shift from the insert above) --------------
$ws.Rows("80:80").Insert()

$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "Macroferia Regional de Talca"
$ws.Range("C80").Value = "Maule"
$ws.Range("D80").Value = 44568
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = 100112001
$ws.Range("G80").Value = "Berenjena"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 200
$ws.Range("K80").Value = 8000
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = 8000
$ws.Range("N80").Value = "$/caja 50 unidades"
$ws.Range("O80").Value = "Región del Maule"
$ws.Range("P80").Value = 160
$ws.Range("Q80").Value = 50
$ws.Range("R80").Value = "Hortaliza"
